$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells (columns F and G) for rows with revised figures
$ws.Cells.Item(244, 6).Value = 5135
$ws.Cells.Item(244, 7).Value = 91
$ws.Cells.Item(246, 6).Value = 2095
$ws.Cells.Item(246, 7).Value = 138
$ws.Cells.Item(250, 6).Value = 17501
$ws.Cells.Item(250, 7).Value = 1147
$ws.Cells.Item(264, 6).Value = 35184
$ws.Cells.Item(264, 7).Value = 800
$ws.Cells.Item(271, 6).Value = 43747
$ws.Cells.Item(271, 7).Value = 1739
$ws.Cells.Item(272, 7).Value = 1637
$ws.Cells.Item(273, 6).Value = 26824
$ws.Cells.Item(276, 6).Value = 11724
$ws.Cells.Item(276, 7).Value = 417
$ws.Cells.Item(278, 6).Value = 30014
$ws.Cells.Item(278, 7).Value = 2072
$ws.Cells.Item(280, 6).Value = 34622
$ws.Cells.Item(280, 7).Value = 2332
$ws.Cells.Item(281, 6).Value = 45522
$ws.Cells.Item(281, 7).Value = 3161
$ws.Cells.Item(282, 6).Value = 47146
$ws.Cells.Item(282, 7).Value = 2842
$ws.Cells.Item(285, 6).Value = 41846
$ws.Cells.Item(285, 7).Value = 3427
$ws.Cells.Item(286, 6).Value = 55084
$ws.Cells.Item(286, 7).Value = 4284
$ws.Cells.Item(287, 6).Value = 58727
$ws.Cells.Item(287, 7).Value = 3715
$ws.Cells.Item(288, 6).Value = 58306
$ws.Cells.Item(288, 7).Value = 4026
$ws.Cells.Item(289, 6).Value = 62806
$ws.Cells.Item(289, 7).Value = 3626
$ws.Cells.Item(290, 6).Value = 17426
$ws.Cells.Item(290, 7).Value = 1034
$ws.Cells.Item(291, 6).Value = 14929
$ws.Cells.Item(291, 7).Value = 485
$ws.Cells.Item(292, 6).Value = 81799
$ws.Cells.Item(293, 6).Value = 81644
$ws.Cells.Item(293, 7).Value = 5743
$ws.Cells.Item(294, 6).Value = 92130
$ws.Cells.Item(294, 7).Value = 4865
$ws.Cells.Item(295, 6).Value = 17199
$ws.Cells.Item(295, 7).Value = 1036
$ws.Cells.Item(298, 6).Value = 3154
$ws.Cells.Item(298, 7).Value = 290
$ws.Cells.Item(300, 6).Value = 71399
$ws.Cells.Item(300, 7).Value = 7033
$ws.Cells.Item(301, 6).Value = 70661
$ws.Cells.Item(301, 7).Value = 5590
$ws.Cells.Item(302, 6).Value = 77043
$ws.Cells.Item(302, 7).Value = 5689
$ws.Cells.Item(304, 6).Value = 5933
$ws.Cells.Item(304, 7).Value = 517
$ws.Cells.Item(305, 6).Value = 3266
$ws.Cells.Item(305, 7).Value = 263
$ws.Cells.Item(306, 7).Value = 7136
$ws.Cells.Item(307, 7).Value = 6418
$ws.Cells.Item(308, 6).Value = 15806
$ws.Cells.Item(308, 7).Value = 1101
$ws.Cells.Item(309, 6).Value = 74236
$ws.Cells.Item(309, 7).Value = 5260
$ws.Cells.Item(310, 6).Value = 75114
$ws.Cells.Item(310, 7).Value = 3912
$ws.Cells.Item(311, 6).Value = 62463
$ws.Cells.Item(313, 6).Value = 71389
$ws.Cells.Item(313, 7).Value = 3208
$ws.Cells.Item(315, 6).Value = 55814
$ws.Cells.Item(315, 7).Value = 2635
$ws.Cells.Item(317, 6).Value = 61972
$ws.Cells.Item(317, 7).Value = 2130
$ws.Cells.Item(318, 6).Value = 49108
$ws.Cells.Item(318, 7).Value = 1193
$ws.Cells.Item(320, 6).Value = 76017
$ws.Cells.Item(320, 7).Value = 3649
$ws.Cells.Item(321, 6).Value = 90514
$ws.Cells.Item(321, 7).Value = 2798
$ws.Cells.Item(322, 6).Value = 106371
$ws.Cells.Item(322, 7).Value = 2294
$ws.Cells.Item(323, 6).Value = 212186
$ws.Cells.Item(323, 7).Value = 3156
$ws.Cells.Item(324, 6).Value = 232920
$ws.Cells.Item(324, 7).Value = 2655
$ws.Cells.Item(325, 6).Value = 752798
$ws.Cells.Item(325, 7).Value = 6356
$ws.Cells.Item(326, 6).Value = 427738
$ws.Cells.Item(326, 7).Value = 3758
$ws.Cells.Item(327, 6).Value = 238950
$ws.Cells.Item(328, 6).Value = 180472
$ws.Cells.Item(329, 6).Value = 88706
$ws.Cells.Item(329, 7).Value = 1792
$ws.Cells.Item(330, 6).Value = 70464
$ws.Cells.Item(330, 7).Value = 1972
$ws.Cells.Item(331, 6).Value = 150049
$ws.Cells.Item(331, 7).Value = 2565
$ws.Cells.Item(332, 6).Value = 424048
$ws.Cells.Item(332, 7).Value = 4126
$ws.Cells.Item(333, 6).Value = 258230
$ws.Cells.Item(333, 7).Value = 2684
$ws.Cells.Item(334, 7).Value = 3377
$ws.Cells.Item(335, 6).Value = 129006
$ws.Cells.Item(335, 7).Value = 2869
$ws.Cells.Item(336, 6).Value = 99374
$ws.Cells.Item(336, 7).Value = 3138
$ws.Cells.Item(337, 6).Value = 100512
$ws.Cells.Item(337, 7).Value = 2839
$ws.Cells.Item(338, 6).Value = 216272
$ws.Cells.Item(338, 7).Value = 3036
$ws.Cells.Item(339, 6).Value = 625181
$ws.Cells.Item(339, 7).Value = 5347
$ws.Cells.Item(340, 6).Value = 370890
$ws.Cells.Item(340, 7).Value = 3137
$ws.Cells.Item(341, 6).Value = 293132
$ws.Cells.Item(341, 7).Value = 3600
$ws.Cells.Item(342, 6).Value = 185085
$ws.Cells.Item(342, 7).Value = 3126
$ws.Cells.Item(343, 6).Value = 121551
$ws.Cells.Item(343, 7).Value = 2716
$ws.Cells.Item(344, 6).Value = 122458
$ws.Cells.Item(344, 7).Value = 2392

# Append new rows 345:347 with full data
$ws.Cells.Item(345, 1).Value = 44239
$ws.Cells.Item(345, 2).Value = 276234
$ws.Cells.Item(345, 3).Value = 13991
$ws.Cells.Item(345, 4).Value = 2330
$ws.Cells.Item(345, 5).Value = 5812
$ws.Cells.Item(345, 6).Value = 250553
$ws.Cells.Item(345, 7).Value = 3087
$ws.Cells.Item(346, 1).Value = 44240
$ws.Cells.Item(346, 2).Value = 277682
$ws.Cells.Item(346, 3).Value = 7597
$ws.Cells.Item(346, 4).Value = 1448
$ws.Cells.Item(346, 5).Value = 5885
$ws.Cells.Item(346, 6).Value = 402206
$ws.Cells.Item(346, 7).Value = 3117
$ws.Cells.Item(347, 1).Value = 44241
$ws.Cells.Item(347, 2).Value = 278254
$ws.Cells.Item(347, 3).Value = 2775
$ws.Cells.Item(347, 4).Value = 572
$ws.Cells.Item(347, 5).Value = 5952
$ws.Cells.Item(347, 6).Value = 249092
$ws.Cells.Item(347, 7).Value = 2392
